# Fruta / hortaliza, semanal
# Insert two new weekly price rows for "Ají" (Inferno, Extra / Primera) at row 415,
# pushing the existing data (rows 415-467) down to rows 417-469.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 415-416; Excel shifts rows 415:467 down to 417:469
# and new blank rows inherit formatting (date style) from the row above.
$ws.Rows("415:416").Insert()

# ---- New row 415 ----
$ws.Range("A415").Value = 10
$ws.Range("B415").Value = "Vega Modelo de Temuco"
$ws.Range("C415").Value = "La Araucanía"
$ws.Range("D415").Value = 44474
$ws.Range("E415").Value = 9
$ws.Range("F415").Value = 100112021
$ws.Range("G415").Value = "Ají"
$ws.Range("H415").Value = "Inferno"
$ws.Range("I415").Value = "Extra"
$ws.Range("J415").Value = 10
$ws.Range("K415").Value = 50000
$ws.Range("L415").Value = 50000
$ws.Range("M415").Value = 50000
$ws.Range("N415").Value = "$/caja 15 kilos"
$ws.Range("O415").Value = "Región de Arica y Parinacota"
$ws.Range("P415").Value = 3333
$ws.Range("Q415").Value = 15
$ws.Range("R415").Value = "Hortaliza"

# ---- New row 416 ----
$ws.Range("A416").Value = 10
$ws.Range("B416").Value = "Vega Modelo de Temuco"
$ws.Range("C416").Value = "La Araucanía"
$ws.Range("D416").Value = 44474
$ws.Range("E416").Value = 9
$ws.Range("F416").Value = 100112021
$ws.Range("G416").Value = "Ají"
$ws.Range("H416").Value = "Inferno"
$ws.Range("I416").Value = "Primera"
$ws.Range("J416").Value = 30
$ws.Range("K416").Value = 45000
$ws.Range("L416").Value = 45000
$ws.Range("M416").Value = 45000
$ws.Range("N416").Value = "$/caja 15 kilos"
$ws.Range("O416").Value = "Región de Arica y Parinacota"
$ws.Range("P416").Value = 3000
$ws.Range("Q416").Value = 15
$ws.Range("R416").Value = "Hortaliza"
